$d = $word.ActiveDocument

# --- Locate the "Example 4:" paragraph and the "135" list-item paragraph ---
# (We find them by content rather than a hard-coded paragraph index so the
#  script is resilient to the document's exact paragraph numbering.)

$exampleParaIndex = 0
$listParaIndex = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Example 4:") {
        $exampleParaIndex = $i
    }
    if ($exampleParaIndex -gt 0 -and $i -eq ($exampleParaIndex + 1) -and $txt -eq "135") {
        $listParaIndex = $i
    }
}

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- 1) Add <w:lastRenderedPageBreak/> before the "Example " run -----------
# Replace the whole "Example 4:" paragraph text with an equivalent set of
# runs, the first of which now carries <w:lastRenderedPageBreak/> ahead of
# its text - mirrors Word re-flowing the break onto this run.
if ($exampleParaIndex -gt 0) {
    $p = $d.Paragraphs($exampleParaIndex)
    $full = $p.Range
    $target = $d.Range($full.Start, $full.Start + 10)   # "Example 4:" (10 chars)

    $xml = "<w:p $wNs>" +
        "<w:r><w:rPr><w:color w:val='7F7F7F' w:themeColor='text1' w:themeTint='80'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:u w:val='thick'/><w:lang w:val='en-US'/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space='preserve'>Example </w:t></w:r>" +
        "<w:r><w:rPr><w:color w:val='7F7F7F' w:themeColor='text1' w:themeTint='80'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:u w:val='thick'/><w:lang w:val='en-US'/></w:rPr><w:t>4</w:t></w:r>" +
        "<w:r><w:rPr><w:color w:val='7F7F7F' w:themeColor='text1' w:themeTint='80'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:u w:val='thick'/><w:lang w:val='en-US'/></w:rPr><w:t>:</w:t></w:r>" +
        "</w:p>"

    $target.InsertXML($xml)
}

# --- 2) Remove <w:lastRenderedPageBreak/> from before the "1" run ----------
# Replace the whole "135" paragraph text with the same two runs, minus the
# page-break marker on the first one.
if ($listParaIndex -gt 0) {
    $p2 = $d.Paragraphs($listParaIndex)
    $full2 = $p2.Range
    $target2 = $d.Range($full2.Start, $full2.Start + 3)  # "135" (3 chars)

    $xml2 = "<w:p $wNs>" +
        "<w:r><w:rPr><w:color w:val='00B050'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:val='en-US'/></w:rPr><w:t>1</w:t></w:r>" +
        "<w:r><w:rPr><w:color w:val='00B050'/><w:sz w:val='28'/><w:szCs w:val='28'/><w:lang w:val='en-US'/></w:rPr><w:t>35</w:t></w:r>" +
        "</w:p>"

    $target2.InsertXML($xml2)
}

Write-Output "exampleParaIndex=$exampleParaIndex listParaIndex=$listParaIndex"
